$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.818.82"
$ws.Range("E2").Value = '  -1.65%  '
$ws.Range("D3").Value = "'1.892.70"
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = "'0.7800"
$ws.Range("E5").Value = '  -3.67%  '
$ws.Range("D6").Value = "'244.17"
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").Value = "'0.9997"
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").Value = "'0.3130"
$ws.Range("E8").Value = '  -3.46%  '
$ws.Range("D9").Value = "'25.42"
$ws.Range("E9").Value = '  -6.50%  '
$ws.Range("D10").Value = "'0.07197"
$ws.Range("E10").Value = '  +1.30%  '
$ws.Range("D11").Value = "'0.08081"
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("D12").Value = "'0.7671"
$ws.Range("E12").Value = '  -2.24%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'5.482"
$ws.Range("E13").Value = '  +1.20%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = "'1.917.66"
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("D15").Value = "'92.46"
$ws.Range("E15").Value = '  -2.57%  '
$ws.Range("D16").Value = "'6.188"
$ws.Range("E16").Value = '  +2.60%  '
$ws.Range("D17").Value = "'29.850.22"
$ws.Range("E17").Value = '  -1.58%  '
$ws.Range("D18").Value = "'13.98"
$ws.Range("E18").Value = '  -2.42%  '
$ws.Range("D19").Value = "'243.85"
$ws.Range("E19").Value = '  -2.39%  '
$ws.Range("D20").Value = "'0.000007771"
$ws.Range("E20").Value = '  -0.69%  '
$ws.Range("D21").Value = "'0.9999"
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("D22").Value = "'8.171"
$ws.Range("E22").Value = '  +3.47%  '
$ws.Range("D23").Value = "'2.148.70"
$ws.Range("E23").Value = '  -1.27%  '
$ws.Range("D24").Value = "'0.9993"
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("D25").Value = "'0.1572"
$ws.Range("E25").Value = '  -3.44%  '
$ws.Range("D26").Value = "'9.408"
$ws.Range("E26").Value = '  -1.17%  '
$ws.Range("D27").Value = "'162.18"
$ws.Range("E27").Value = '  -3.33%  '
$ws.Range("D28").Value = "'18.75"
$ws.Range("E28").Value = '  -1.90%  '
$ws.Range("D29").Value = "'2.052"
$ws.Range("E29").Value = '  -3.59%  '
$ws.Range("D30").Value = "'1.424"
$ws.Range("E30").Value = '  +3.65%  '
$ws.Range("D31").Value = "'1.551"
$ws.Range("E31").Value = '  +0.74%  '
$ws.Range("D32").Value = "'4.483"
$ws.Range("E32").Value = '  +2.75%  '
$ws.Range("D33").Value = "'4.108"
$ws.Range("E33").Value = '  -0.66%  '
$ws.Range("D34").Value = "'0.05534"
$ws.Range("E34").Value = '  -1.97%  '
$ws.Range("D35").Value = "'1.262"
$ws.Range("E35").Value = '  -3.28%  '
$ws.Range("D36").Value = "'0.7482"
$ws.Range("E36").Value = '  +0.84%  '
$ws.Range("D37").Value = "'1.005"
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("D38").Value = "'2.627"
$ws.Range("E38").Value = '  -3.40%  '
$ws.Range("D39").Value = "'0.01920"
$ws.Range("E39").Value = '  -1.61%  '
$ws.Range("D40").Value = "'2.778"
$ws.Range("E40").Value = '  -1.53%  '
$ws.Range("D41").Value = "'1.141.87"
$ws.Range("E41").Value = '  +10.05%  '
$ws.Range("D42").Value = "'73.83"
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").Value = "'0.4427"
$ws.Range("E43").Value = '  -1.21%  '
$ws.Range("D44").Value = "'5.906"
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("D45").Value = "'0.8521"
$ws.Range("E45").Value = '  -0.44%  '
$ws.Range("D46").Value = "'104.02"
$ws.Range("E46").Value = '  +0.87%  '
$ws.Range("D47").Value = "'0.9996"
$ws.Range("E47").Value = '  -0.25%  '
$ws.Range("D48").Value = "'1.896"
$ws.Range("E48").Value = '  -2.05%  '
$ws.Range("D49").Value = "'9.925"
$ws.Range("E49").Value = '  -0.66%  '
$ws.Range("D50").Value = "'3.045"
$ws.Range("E50").Value = '  +11.43%  '
$ws.Range("D51").Value = "'7.470"
$ws.Range("E51").Value = '  -2.32%  '
